$wb = $excel.ActiveWorkbook

# --- Sheet1 ("超大字符集"): append 4 new rows (146-149) describing the
# plane02 "min-l"/"min-r" variants, mirroring the existing row layout.
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A146").Value = "up-min-l-plane02"
$ws1.Range("B146").Value = "unicode"
$ws1.Range("C146").Value = "ZhongHuaSongPlane02-Regular.ttf"

$ws1.Range("A147").Value = "ut-min-l-plane02"
$ws1.Range("B147").Value = "unicode"
$ws1.Range("C147").Value = "ZhongHuaSongPlane02-Regular.ttf"
$ws1.Range("D147").Value = "-w 1"

$ws1.Range("A148").Value = "up-min-r-plane02"
$ws1.Range("B148").Value = "unicode"
$ws1.Range("C148").Value = "ZhongHuaSongPlane02-Regular.ttf"

$ws1.Range("A149").Value = "ut-min-r-plane02"
$ws1.Range("B149").Value = "unicode"
$ws1.Range("C149").Value = "ZhongHuaSongPlane02-Regular.ttf"
$ws1.Range("D149").Value = "-w 1"

# --- Sheet2 ("URO-Utmin10"): selection moves to B50. Do this first, while
# sheet2 is still the active sheet, so selecting the range doesn't flip the
# workbook's active tab back afterwards.
$ws2 = $wb.Worksheets.Item(2)
[void]$ws2.Range("B50").Select()

# --- View state: sheet1 ("超大字符集") becomes the active tab/sheet
# (workbook's activeTab drops the explicit "1" it had before, sheet1 gains
# tabSelected, sheet2 loses it), with the on-sheet selection moved to A150.
[void]$ws1.Activate()
[void]$ws1.Range("A150").Select()
